$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

${ws}.Range('A2').Value = 'KEsQQ150'
${ws}.Range("B2").Value = 23102849
${ws}.Range('C2').Value = 'zouzcrj59'
${ws}.Range('D2').Value = 'AZ#p$2r7'
${ws}.Range('F2').Value = 'sctXiARj'
${ws}.Range('G2').Value = 'uvnu'
${ws}.Range('A3').Value = 'JgMgt863'
${ws}.Range("B3").Value = 23102848
${ws}.Range('C3').Value = 'uolunsd76'
${ws}.Range('D3').Value = 'WrV8#%4b'
${ws}.Range('F3').Value = 'fJQvOcrR'
${ws}.Range('G3').Value = 'OMOT'
${ws}.Range('A4').Value = 'vrGvq212'
${ws}.Range("B4").Value = 23102847
${ws}.Range('C4').Value = 'gcsvsqv20'
${ws}.Range('D4').Value = 'S98H!xh&'
${ws}.Range('F4').Value = 'BnWFmmez'
${ws}.Range('G4').Value = 'mcUh'
${ws}.Range('A5').Value = 'uIzKV818'
${ws}.Range("B5").Value = 23102846
${ws}.Range('C5').Value = 'tqaxvul92'
${ws}.Range('D5').Value = 'X6E!%m8n'
${ws}.Range('F5').Value = 'gJofDxwm'
${ws}.Range('G5').Value = 'wcDK'
${ws}.Range('A6').Value = 'DsTQc634'
${ws}.Range("B6").Value = 23102845
${ws}.Range('C6').Value = 'lbygnmo24'
${ws}.Range('D6').Value = 't$%XG97x'
${ws}.Range('F6').Value = 'MjtfvdFs'
${ws}.Range('G6').Value = 'QNIz'
${ws}.Range('A7').Value = 'ZkWVx405'
${ws}.Range("B7").Value = 23102844
${ws}.Range('C7').Value = 'egrtzvc51'
${ws}.Range('D7').Value = 'ux#G!26F'
${ws}.Range('F7').Value = 'Kighziev'
${ws}.Range('G7').Value = 'pjfx'
${ws}.Range('A8').Value = 'yglHs934'
${ws}.Range("B8").Value = 23102843
${ws}.Range('C8').Value = 'xihqopy62'
${ws}.Range('D8').Value = 'H6n$A!4f'
${ws}.Range('F8').Value = 'WCJqNpMp'
${ws}.Range('G8').Value = 'LuOL'
${ws}.Range('A9').Value = 'xdHHH982'
${ws}.Range("B9").Value = 23102842
${ws}.Range('C9').Value = 'cptyjlt63'
${ws}.Range('D9').Value = 'P26!T#gv'
${ws}.Range('F9').Value = 'CwFRsviz'
${ws}.Range('G9').Value = 'wHFx'
${ws}.Range('A10').Value = 'LpPKY107'
${ws}.Range("B10").Value = 23102841
${ws}.Range('C10').Value = 'wennkzu38'
${ws}.Range('D10').Value = 'q6CVy%#3'
${ws}.Range('F10').Value = 'OqWAaXkS'
${ws}.Range('G10').Value = 'QacE'
${ws}.Range('A11').Value = 'deVvo766'
${ws}.Range("B11").Value = 23102840
${ws}.Range('C11').Value = 'zdbvgfx20'
${ws}.Range('D11').Value = 'q#7vXM$6'
${ws}.Range('F11').Value = 'jSeINGcv'
${ws}.Range('G11').Value = 'ZVgI'
${ws}.Range('A12').Value = 'rvJaC610'
${ws}.Range("B12").Value = 23102839
${ws}.Range('C12').Value = 'cerleoe53'
${ws}.Range('D12').Value = 'XV4hy$5#'
${ws}.Range('F12').Value = 'lYSHxzlz'
${ws}.Range('G12').Value = 'TERW'
${ws}.Range('A13').Value = 'tMKFb811'
${ws}.Range("B13").Value = 23102838
${ws}.Range('C13').Value = 'rzqbkvx16'
${ws}.Range('D13').Value = 'tN5!Ym7&'
${ws}.Range('F13').Value = 'FiPkRntl'
${ws}.Range('G13').Value = 'vMda'
${ws}.Range('A14').Value = 'jmyfI815'
${ws}.Range("B14").Value = 23102837
${ws}.Range('C14').Value = 'dkcvjrg64'
${ws}.Range('D14').Value = 'g7bWV%!6'
${ws}.Range('F14').Value = 'yQOHRmjK'
${ws}.Range('G14').Value = 'mZMB'
${ws}.Range('A15').Value = 'gOIzY792'
${ws}.Range("B15").Value = 23102836
${ws}.Range('C15').Value = 'ybznags33'
${ws}.Range('D15').Value = 'nW36k!J&'
${ws}.Range('F15').Value = 'tdrWiRvw'
${ws}.Range('G15').Value = 'AovN'
